$wb = $excel.ActiveWorkbook

# --- Sheet "2018" : row 2 (2018 Q2) and row 3 (2018 Q3) ---
$ws2018 = $wb.Worksheets.Item("2018")

# Row 2
$ws2018.Range("C2").Value = 465890
$ws2018.Range("D2").Value = 410500
$ws2018.Range("F2").Value = 454450
$ws2018.Range("G2").Value = 382500
$ws2018.Range("I2").Value = 10
$ws2018.Range("J2").Value = 456.11

# Row 3
$ws2018.Range("C3").Value = 478308.25
$ws2018.Range("D3").Value = 470666.6666666667
$ws2018.Range("F3").Value = 499000
$ws2018.Range("G3").Value = 487500
$ws2018.Range("I3").Value = 12
$ws2018.Range("J3").Value = 362.05

# --- Sheet "2019" : row 3 (2018 Q4) ---
$ws2019 = $wb.Worksheets.Item("2019")

$ws2019.Range("C3").Value = 462560
$ws2019.Range("D3").Value = 433580
$ws2019.Range("G3").Value = 469900
$ws2019.Range("I3").Value = 10
$ws2019.Range("J3").Value = 481.76
